$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the three tables (slides 14, 15, 16) from the default
#    "Table_0" style to the built-in "{A499E4C2-9741-4252-924F-15AE7BFFE762}"
#    table style.
# ---------------------------------------------------------------------------
$newTableStyle = "{A499E4C2-9741-4252-924F-15AE7BFFE762}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the presentation's theme colour scheme back to the stock
#    "Office" palette (it currently carries the "Integral" / Red Violet
#    palette). All slides share a single master/theme, so touching the
#    scheme from any one slide updates it for the whole deck.
# ---------------------------------------------------------------------------
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$scheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
